# encargados.xlsx: emails are now managed solely from config.json, so the
# "mail" column (with its mailto hyperlinks) is no longer needed in the
# spreadsheet. Remove it and the hyperlink styling that went with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the mailto: hyperlinks that lived on column J ("mail") before the
# column itself is removed.
$ws.Hyperlinks.Delete()

# Remove the whole "mail" column (J) - "seccion/es" (formerly column K)
# shifts left into J.
$ws.Columns("J").Delete()

# The hyperlink cell style ("Hipervínculo") is no longer used anywhere in
# the sheet now that the mail column/hyperlinks are gone - drop it so it
# doesn't linger in the style table.
$wb.Styles.Item("Hipervínculo").Delete()

# Match the author's last on-screen selection after the edit.
$ws.Range("L10").Select()
